# Dataframe ST.xlsx - apply "Add files via upload" edit
#  - Sheet3: re-key the lookup block (B2:B18) to match Sheet1's B2:B18 order,
#    refresh the raw data block (B20:B36) with new figures, and turn C2:C18
#    into IFERROR(VLOOKUP(...,A20:B36,2,),0) formulas (C2 alone, C3:C18 shared).
#  - Sheet1: add a new date column CD ("30-oct") that mirrors CB/CC with the
#    same VLOOKUP-into-Sheet3 formula, per row.
#  - Restore the selections Excel left behind (Sheet1 CB7 active, Sheet3 C4
#    active, Sheet1 still the visible/selected tab).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# 1) Sheet3!B2:B18 - same descriptions as Sheet1!B2:B18, in Sheet1's row order
# ---------------------------------------------------------------------------
$descByRow = @{
    2  = "LAYS CLASICAS 40GX68"
    3  = "LAYS CLASICAS 94GRX25"
    4  = "LAYS CLASICAS 145GRX18"
    5  = "LAYS CLASICAS 249GRX14"
    6  = "DORITOS QUESO 40GX58X1 CH"
    7  = "DORITOS QUESO 85GX26"
    8  = "DORITOS QUESO 140GX19"
    9  = "PEHUAMAR PAPA LISA 520GX9"
    10 = "PEHUAMAR ACANALADA 520GX9"
    11 = "PEHUAMAR MAICITOS 285GX10"
    12 = "3D QUESO 92GX27"
    13 = "CHEETOS 94GRX24"
    14 = "QUAKER AVENA INSTANT FORTIF 18X280G"
    15 = "LAYS QSO Y CEBOLLA 34GX72"
    16 = "LAYS CEBOLLA CARAMELIZADA 85GX25"
    17 = "LAYS ONDAS FH 30GX72"
    18 = "LAYS ONDAS FH 70GX28"
}
for ($r = 2; $r -le 18; $r++) {
    $ws3.Cells.Item($r, 2).Value = $descByRow[$r]
}

# ---------------------------------------------------------------------------
# 2) Sheet3!B20:B36 - refreshed raw figures behind the lookup table
# ---------------------------------------------------------------------------
$rawByRow = @{
    20 = 14.969170967259931
    21 = 3.1391388891607512
    22 = 12.7582050172142
    23 = 0
    24 = 0
    25 = 0
    26 = 2.3961915521389523
    27 = 8.2812015243249526
    28 = 4.9276109148536538
    29 = 3.888431220736023
    30 = 0.21330509147019314
    31 = 3.0901251487345984
    32 = 5.6742696089004045
    33 = 6.2279493659965617
    34 = 2.2047859666365399
    35 = 8.6610022930455202
    36 = 42.153156108638647
}
for ($r = 20; $r -le 36; $r++) {
    $ws3.Cells.Item($r, 2).Value = $rawByRow[$r]
}

# ---------------------------------------------------------------------------
# 3) Sheet3!C2:C18 - IFERROR(VLOOKUP(...),0) against the A20:B36 block
#    (C2 on its own, C3:C18 entered together so Excel shares the formula)
# ---------------------------------------------------------------------------
$ws3.Range("C2").Formula = '=IFERROR(VLOOKUP(B2,A20:B36,2,),0)'
$ws3.Range("C3:C18").Formula = '=IFERROR(VLOOKUP(B3,A21:B37,2,),0)'

# ---------------------------------------------------------------------------
# 4) Sheet1 - new column CD ("30-oct"), mirroring CB/CC
# ---------------------------------------------------------------------------
$ws1.Range("CD1").Value = "30-oct"

for ($r = 2; $r -le 18; $r++) {
    $cell = $ws1.Cells.Item($r, 82)
    $cell.NumberFormat = "0"
    $cell.Formula = '=VLOOKUP($B' + $r + ',Sheet3!$B$1:$C$18,2,)'
}

# ---------------------------------------------------------------------------
# 5) Selections: Sheet3 -> C4, then back to Sheet1 -> CB7 so Sheet1 stays the
#    active/visible tab (matches the saved view state in the workbook).
# ---------------------------------------------------------------------------
$ws3.Range("C4").Select()
$ws1.Range("CB7").Select()
